$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 97-99 near the "Other" section ---
# Old row 98 (power supply / 1 / 7) moves down to row 99
$ws.Range("B99").Value = $ws.Range("B98").Value2
$ws.Range("D99").Value = $ws.Range("D98").Value2
$ws.Range("E99").Value = $ws.Range("E98").Value2

# Old row 97 (30 guage wire / approx 36') moves down to row 98
$ws.Range("B98").Value = $ws.Range("B97").Value2
$ws.Range("F98").Value = $ws.Range("F97").Value2
$ws.Range("D98").Clear() | Out-Null
$ws.Range("E98").Clear() | Out-Null

# New row 97 content: the 6mm diametric magnet included with the as5600 sensors
$ws.Range("B97").Value = "6mm diametric magnet"
$ws.Range("D97").Value = 3
$ws.Range("F97").Value = "included with as5600 sensors"

# --- Row 61: bump quantity and change the pcb note to "pcb, sensor" ---
$ws.Range("D61").Value = 9
$ws.Range("F61").Value = "pcb, sensor"

# --- Sheet view state (scroll position + active selection) ---
$window = $ws.Application.ActiveWindow
$window.ScrollRow = 31
$window.ScrollColumn = 1
$ws.Range("F62").Select() | Out-Null
